$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append two new data rows (15 and 16) below the existing table, inheriting
# the formatting of the row above (date style for column A, time style for
# column B) by inserting rows rather than just writing into blank cells.

$ws.Rows("15").Insert(-4121) | Out-Null   # xlShiftDown
$ws.Range("A15").Value = 41568            # 21/10/2013
$ws.Range("B15").Value = 0.024305555555555556

$ws.Rows("16").Insert(-4121) | Out-Null   # xlShiftDown
$ws.Range("A16").Value = 41569            # 22/10/2013
$ws.Range("B16").Value = 0.125

# Move the active selection to the next empty cell below the table, matching
# the post-edit state of the workbook.
$ws.Range("B17").Select() | Out-Null
